$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated feature importance table (features reordered/renamed, importances recalculated)
$features = @(
    'decisão proferida',
    'aumento',
    'xliii',
    'concessão',
    'ordem prisão',
    'justiça indeferiu',
    'julgado turma',
    'arquivado',
    'enunciado súmula',
    'deferida',
    'senha relatório',
    'concessão ordem',
    'outro motivo',
    'liberdade provisória',
    'assessoria',
    'implicou deferimento',
    'deferido',
    'liminar espécie',
    'stj indeferiu',
    'liminar hc',
    'deferi',
    'relatório http',
    'restritiva direitos',
    'manifesta ilegalidade',
    'opina deferimento',
    'restritiva',
    'informado',
    'http sob',
    'inconstitucional',
    'provisória formulado',
    'cautelar pois',
    'proferida ministro',
    'preventiva fundamentos',
    'tutela cautelar',
    'fundamentos insubsistência',
    'liminar suspender',
    'lxi',
    'opinou concessão',
    'arquivado definitivo',
    'república concessão',
    'gabinete prestou',
    'formalizado ato',
    'deferimento liminar',
    'deferida assessoria',
    'deferi pedido',
    'sobrestamento',
    'código senha',
    'contornos impetração',
    'desta suprema',
    'definitivo após',
    'efêmero',
    'resumida prisão',
    'senha primeira',
    'benefício liberdade',
    'ficou',
    'flagrante tráfico',
    'sob código',
    'precário',
    'pertence grifei',
    'precário efêmero',
    'suspender efeitos',
    'sistema jurídico',
    'senha',
    'princípio liberdade',
    'rtj rel',
    'processo formalizado',
    'revelou contornos',
    'punir',
    'suspender',
    'paulo indeferiu',
    'resumida',
    'análise pedido',
    'mendes hc',
    'liminar deferida',
    'art lvii',
    'art lxi',
    'assessoria prestou',
    'assim resumida',
    'assim revelou',
    'aurélio decisão',
    'campo precário',
    'condenação penal',
    'constitucional cf',
    'defiro',
    'deserção',
    'efeitos ordem',
    'eis informado',
    'espécie ficou',
    'fiança',
    'ficou assim',
    'http',
    'impetração eis',
    'inciso xliii',
    'indefiro pedido',
    'infligir',
    'infligir punição',
    'informado análise',
    'instituto prisão',
    'irrecorrível',
    'indiciado réu'
)

$importances = @(
    0.05,
    0.05,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0
)

for ($i = 0; $i -lt $features.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $features[$i]
    $ws.Cells.Item($row, 2).Value = $importances[$i]
}
